# Remove the first two data rows (2008年 and 2009年) from the sheet.
# This shifts the remaining rows (2010年, 2011年) up, so they become
# rows 2 and 3, and the sheet's used range shrinks from A1:AT5 to A1:AT3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
